$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '26.615.52'
Set-TextValue 'E2' '  -0.01%  '
Set-TextValue 'D3' '1.596.00'
Set-TextValue 'E3' '  +0.51%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '211.54'
Set-TextValue 'E5' '  +0.29%  '
Set-TextValue 'D6' '0.516'
Set-TextValue 'E6' '  +1.22%  '
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'E8' '  +0.35%  '
Set-TextValue 'E9' '  +0.02%  '
Set-TextValue 'D10' '19.50'
Set-TextValue 'E10' '  -0.42%  '
Set-TextValue 'E11' '  +0.33%  '
Set-TextValue 'D12' '1.820.41'
Set-TextValue 'E12' '  +0.53%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.601.24'
Set-TextValue 'E13' '  +0.67%  '
Set-TextValue 'B14' 'Polkadot'
Set-TextValue 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D14' '4.03'
Set-TextValue 'E14' '  +0.12%  '
Set-TextValue 'E15' '  -0.01%  '
Set-TextValue 'D16' '64.50'
Set-TextValue 'E16' '  -0.32%  '
Set-TextValue 'D17' '26.606.23'
Set-TextValue 'E17' '  +0.02%  '
Set-TextValue 'D18' '0.0₃0731'
Set-TextValue 'E18' '  +0.51%  '
Set-TextValue 'D19' '208.78'
Set-TextValue 'E19' '  +0.35%  '
Set-TextValue 'E20' '  -0.04%  '
Set-TextValue 'D21' '6.94'
Set-TextValue 'E21' '  +3.23%  '
Set-TextValue 'D22' '4.26'
Set-TextValue 'E22' '  +0.43%  '
Set-TextValue 'E23' '  -1.58%  '
Set-TextValue 'D24' '8.90'
Set-TextValue 'D25' '145.33'
Set-TextValue 'E25' '  -1.00%  '
Set-TextValue 'E26' '  -0.03%  '
Set-TextValue 'D27' '7.13'
Set-TextValue 'E27' '  -1.35%  '
Set-TextValue 'E28' '  +0.89%  '
Set-TextValue 'D29' '15.26'
Set-TextValue 'E29' '  -0.16%  '
Set-TextValue 'D30' '0.0507'
Set-TextValue 'E30' '  +0.03%  '
Set-TextValue 'E31' '  +0.61%  '
Set-TextValue 'E32' '  +0.30%  '
Set-TextValue 'E33' '  -1.61%  '
Set-TextValue 'E34' '  +1.03%  '
Set-TextValue 'D35' '1.281.61'
Set-TextValue 'E35' '  -1.87%  '
Set-TextValue 'E36' '  +1.15%  '
Set-TextValue 'E37' '  +1.09%  '
Set-TextValue 'D39' '0.843'
Set-TextValue 'E39' '  +1.87%  '
Set-TextValue 'E40' '  -0.05%  '
Set-TextValue 'E41' '  +1.86%  '
Set-TextValue 'B42' 'Aave'
Set-TextValue 'C42' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '64.39'
Set-TextValue 'E42' '  +2.85%  '
Set-TextValue 'B43' 'TrustWalletToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D43' '0.785'
Set-TextValue 'E43' '  -0.95%  '
Set-TextValue 'E44' '  +1.23%  '
Set-TextValue 'D45' '1.733.09'
Set-TextValue 'E45' '  +0.52%  '
Set-TextValue 'D46' '0.908'
Set-TextValue 'E46' '  +8.29%  '
Set-TextValue 'D47' '89.67'
Set-TextValue 'E47' '  +0.20%  '
Set-TextValue 'E48' '  -0.15%  '
Set-TextValue 'E49' '  +5.18%  '
Set-TextValue 'D50' '0.0506'
Set-TextValue 'E50' '  +0.35%  '
Set-TextValue 'D51' '7.47'
Set-TextValue 'E51' '  -0.73%  '
